$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.30899965763092
$ws.Range("B1").Value = 3.277969121932983
$ws.Range("C1").Value = 5.700580596923828
$ws.Range("D1").Value = 1.728224158287048
$ws.Range("E1").Value = 1.012715339660645
